$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '51.903.40'
$c.ClearFormats()
$ws.Range('E2').Value = '  +0.30%  '

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.815.38'
$c.ClearFormats()
$ws.Range('E3').Value = '  +1.37%  '

$ws.Range('E4').Value = '  -0.08%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '356.83'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.14%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '110.46'
$c.ClearFormats()
$ws.Range('E6').Value = '  +1.48%  '

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.557'
$c.ClearFormats()
$ws.Range('E7').Value = '  +0.46%  '

$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E8').Value = '  -0.06%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.636'
$c.ClearFormats()
$ws.Range('E9').Value = '  +8.72%  '

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '40.24'
$c.ClearFormats()
$ws.Range('E10').Value = '  +1.63%  '

$ws.Range('E11').Value = '  +0.06%  '

$ws.Range('E12').Value = '  -0.54%  '

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '20.06'
$c.ClearFormats()
$ws.Range('E13').Value = '  +3.27%  '

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '7.82'
$c.ClearFormats()
$ws.Range('E14').Value = '  +2.94%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '3.253.76'
$c.ClearFormats()
$ws.Range('E15').Value = '  +1.28%  '

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '2.812.73'
$c.ClearFormats()
$ws.Range('E16').Value = '  +0.91%  '

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.945'
$c.ClearFormats()
$ws.Range('E17').Value = '  +1.38%  '

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '51.838.55'
$c.ClearFormats()
$ws.Range('E18').Value = '  +0.29%  '

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '7.69'
$c.ClearFormats()
$ws.Range('E19').Value = '  +3.42%  '

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '3.20'
$c.ClearFormats()
$ws.Range('E20').Value = '  +3.57%  '

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '13.72'
$c.ClearFormats()
$ws.Range('E21').Value = '  +4.62%  '

$ws.Range('E22').Value = '  +1.30%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '70.48'
$c.ClearFormats()
$ws.Range('E23').Value = '  +0.45%  '

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '268.64'
$c.ClearFormats()
$ws.Range('E24').Value = '  +0.03%  '

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.77'
$c.ClearFormats()
$ws.Range('E25').Value = '  +1.23%  '

$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '26.21'
$c.ClearFormats()
$ws.Range('E26').Value = '  -0.62%  '

$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range('E27').Value = '  +0.02%  '

$ws.Range('E28').Value = '  -0.86%  '

$ws.Range('E29').Value = '  +2.42%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '37.97'
$c.ClearFormats()
$ws.Range('E30').Value = '  +9.42%  '

$ws.Range('E31').Value = '  -1.84%  '

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '6.21'
$c.ClearFormats()
$ws.Range('E32').Value = '  -0.22%  '

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '52.23'
$c.ClearFormats()
$ws.Range('E33').Value = '  +1.25%  '

$ws.Range('E34').Value = '  +10.92%  '

$ws.Range('E35').Value = '  -1.25%  '

$ws.Range('E36').Value = '  +3.87%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E37').Value = '  -0.07%  '

$ws.Range('E38').Value = '  +1.21%  '

$ws.Range('E39').Value = '  +3.10%  '

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '3.15'
$c.ClearFormats()
$ws.Range('E40').Value = '  +0.60%  '

$ws.Range('E41').Value = '  +1.21%  '

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '2.50'
$c.ClearFormats()
$ws.Range('E42').Value = '  -1.59%  '

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '120.19'
$c.ClearFormats()
$ws.Range('E43').Value = '  +0.82%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '22.10'
$c.ClearFormats()
$ws.Range('E44').Value = '  +1.87%  '

$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.ClearFormats()
$ws.Range('E45').Value = '  -0.61%  '

$ws.Range('E46').Value = '  +8.21%  '

$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '3.39'
$c.ClearFormats()
$ws.Range('E47').Value = '  +4.04%  '

$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '2.109.34'
$c.ClearFormats()
$ws.Range('E48').Value = '  +1.28%  '

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.931'
$c.ClearFormats()
$ws.Range('E49').Value = '  -1.10%  '

$ws.Range('E50').Value = '  +9.68%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '5.46'
$c.ClearFormats()
$ws.Range('E51').Value = '  -1.57%  '
